# E2e test cases added
$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item("BookNow").Name = "HomePage"
$wb.Worksheets.Item("Find&Reserve").Name = "SingnInPage"

# Update cell contents on the (now renamed) HomePage sheet
$ws = $wb.Worksheets.Item("HomePage")
$ws.Range("E2").Value = "2 Rooms"
$ws.Range("E3").Value = "2 Rooms"

# Move the active selection
$ws.Range("D12").Select()
